$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "61.201.67"
$ws.Range("E2").Value = "  -1.93%  "
$ws.Range("D3").Value = "2.981.84"
$ws.Range("E3").Value = "  -0.73%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.12%  "
Set-TextValue $ws.Range("D5") "595.36"
$ws.Range("E5").Value = "  +1.60%  "
Set-TextValue $ws.Range("D6") "143.18"
$ws.Range("E6").Value = "  -1.93%  "
Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "2.979.88"
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue $ws.Range("D9") "0.513"
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("E10").Value = "  -0.76%  "
Set-TextValue $ws.Range("D11") "6.01"
$ws.Range("E11").Value = "  +4.27%  "
Set-TextValue $ws.Range("D12") "0.453"
$ws.Range("E12").Value = "  +2.90%  "
Set-TextValue $ws.Range("D13") "0.0000226"
$ws.Range("E13").Value = "  -0.79%  "
Set-TextValue $ws.Range("D14") "34.15"
$ws.Range("E14").Value = "  -1.04%  "
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").Value = "3.474.03"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").Value = "61.134.59"
$ws.Range("E17").Value = "  -1.98%  "
Set-TextValue $ws.Range("D18") "6.86"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("D19").Value = "2.977.09"
$ws.Range("E19").Value = "  -0.83%  "
Set-TextValue $ws.Range("D20") "448.95"
$ws.Range("E20").Value = "  -1.35%  "
Set-TextValue $ws.Range("D21") "14.01"
$ws.Range("E21").Value = "  +1.42%  "
Set-TextValue $ws.Range("D22") "0.682"
$ws.Range("E22").Value = "  +0.47%  "
Set-TextValue $ws.Range("D23") "7.29"
$ws.Range("E23").Value = "  -1.37%  "
Set-TextValue $ws.Range("D24") "81.77"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D25") "2.17"
$ws.Range("E25").Value = "  -3.99%  "
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D26") "10.44"
$ws.Range("E26").Value = "  +3.36%  "
Set-TextValue $ws.Range("D27") "11.94"
$ws.Range("E27").Value = "  -2.33%  "
$ws.Range("E28").Value = "  -0.02%  "
Set-TextValue $ws.Range("D29") "2.68"
$ws.Range("E29").Value = "  +2.30%  "
Set-TextValue $ws.Range("D30") "0.999"
$ws.Range("E30").Value = "  -0.07%  "
Set-TextValue $ws.Range("D31") "7.20"
$ws.Range("E31").Value = "  +0.20%  "
Set-TextValue $ws.Range("D32") "2.05"
$ws.Range("E32").Value = "  -2.15%  "
Set-TextValue $ws.Range("D33") "27.17"
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").Value = "0.0₃0812"
$ws.Range("E35").Value = "  +3.00%  "
$ws.Range("E36").Value = "  -1.33%  "
Set-TextValue $ws.Range("D37") "5.77"
$ws.Range("E37").Value = "  +0.86%  "
Set-TextValue $ws.Range("D38") "50.09"
$ws.Range("E38").Value = "  +0.14%  "
Set-TextValue $ws.Range("D39") "2.05"
$ws.Range("E39").Value = "  -2.43%  "
Set-TextValue $ws.Range("D40") "8.94"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("E41").Value = "  +5.86%  "
Set-TextValue $ws.Range("D42") "2.83"
$ws.Range("E42").Value = "  -2.80%  "
Set-TextValue $ws.Range("D43") "384.42"
$ws.Range("E43").Value = "  -0.16%  "
Set-TextValue $ws.Range("D44") "0.269"
$ws.Range("E44").Value = "  -1.33%  "
Set-TextValue $ws.Range("D45") "0.0350"
$ws.Range("E45").Value = "  -0.69%  "
Set-TextValue $ws.Range("D46") "38.58"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "2.696.33"
Set-TextValue $ws.Range("D48") "130.49"
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("E50").Value = "  -0.48%  "
Set-TextValue $ws.Range("D51") "2.13"
$ws.Range("E51").Value = "  -0.67%  "
